$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the existing last row (row 38): the "added" date moves back a day ---
$ws.Range("C38").Value = 44383

# --- Append the new sample as row 39 ---
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "msg"

# C39 needs the same date style as the rest of column C (s="1", numFmtId 14).
# Copy the format from C38 (a cell that already has that style) before setting
# the new date value, so we don't introduce a brand-new style entry.
$ws.Range("C38").Copy()
$ws.Range("C39").PasteSpecial(-4122)
$ws.Range("C39").Value = 44384

$ws.Range("D39").Value = "MCAST"
$ws.Range("E39").Value = "lookalike"
$ws.Range("F39").Value = "delivery"
$ws.Range("G39").Value = "mt"
$ws.Range("H39").Value = "no"
$ws.Range("I39").Value = "payment request for postal delivery"
$ws.Range("J39").Value = "DHL"

# --- Re-apply the autofilter / filter-database range to cover the new row ---
[void]($ws.AutoFilterMode = $false)
[void]$ws.Range("A1:K39").AutoFilter()

$fdName = $wb.Names.Item(1)
$fdName.RefersTo = "=Sheet1!`$A`$1:`$K`$39"

# --- Restore the selected cell shown in the sheet view ---
[void]$ws.Range("K31").Select()
